$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 8.938819938647542; C = 5.600988739745013; D = 4.572609394117377; E = 16.48800361265332; F = 22.86306017478281; H = 7.344005520526261; I = 19.31656217468171; K = 8.736392222137528; N = 17.52924621832931; O = 20.43177991603501 }
    3 = @{ B = 8.595745181530914; C = 5.344437272379126; D = 4.515740681239259; E = 15.55441097959246; F = 22.85798319484641; H = 7.344005520526261; I = 19.39331494579942; K = 8.49156782699146; N = 17.58768085849675; O = 20.48199027864781 }
    4 = @{ B = 8.379316793265259; C = 5.179321475052843; D = 4.479922716352814; E = 14.95632413580999; F = 22.86150698468045; H = 7.344005520526261; I = 19.44440815652837; K = 8.339027531667458; N = 17.62520077763389; O = 20.51773012443044 }
    5 = @{ B = 8.289805315195817; C = 5.110174352179126; D = 4.46510797389413; E = 14.70661821184783; F = 22.86461198251146; H = 7.344005520526261; I = 19.46622464860719; K = 8.276404596909034; N = 17.6409043639107; O = 20.53352538254939 }
    6 = @{ B = 8.274866895163072; C = 5.098581682426955; D = 4.462635023427371; E = 14.66480244666628; F = 22.8652283180901; H = 7.344005520526261; I = 19.46990733749211; K = 8.265981147734601; N = 17.6435369706467; O = 20.53622240527008 }
    7 = @{ B = 8.378114741901937; C = 5.178396398895961; D = 4.479723794386923; E = 14.95298033449857; F = 22.86154210416057; H = 7.344005520526261; I = 19.4446983527408; K = 8.338184714561979; N = 17.62541088403106; O = 20.51793816586434 }
    8 = @{ B = 8.821805501660798; C = 5.514131926118308; D = 4.553192520359587; E = 16.17140624897705; F = 22.85993109341511; H = 7.344005520526261; I = 19.34220212193973; K = 8.652492726037789; N = 17.54905472874943; O = 20.4480710686587 }
    9 = @{ B = 9.640695799583842; C = 6.11053979213652; D = 4.689791849368132; E = 18.43668239167303; F = 22.90945244221208; H = 7.344005520526261; I = 19.17276468284195; K = 9.2473432046315; N = 17.41227847914674; O = 20.35018428811927 }
    10 = @{ B = 10.20495719013449; C = 6.50908653395937; D = 4.785174179689933; E = 20.07305316876846; F = 22.97785299817646; H = 7.344005520526261; I = 19.06763133117407; K = 9.66647281377713; N = 17.31960138553531; O = 20.3023204884332 }
    11 = @{ B = 10.45248734535766; C = 6.681506434957936; D = 4.827401096023387; E = 20.77527044936654; F = 23.01587354065688; H = 7.344005520526261; I = 19.02403025387794; K = 9.852373386948601; N = 17.27911795372293; O = 20.28580655244902 }
    12 = @{ B = 10.54483155020262; C = 6.745502237611032; D = 4.843217351457295; E = 21.03515283316822; F = 23.03125743144054; H = 7.344005520526261; I = 19.00812924808023; K = 9.922019169569534; N = 17.26402753495288; O = 20.28031199204038 }
    13 = @{ B = 10.52500646403383; C = 6.731777491317597; D = 4.839818890874803; E = 20.97945026964604; F = 23.02790048672494; H = 7.344005520526261; I = 19.0115266579373; K = 9.907054073074585; N = 17.26726688153129; O = 20.28146155369965 }
    14 = @{ B = 10.46011280418878; C = 6.686797474023955; D = 4.828705837955241; E = 20.79677183394257; F = 23.01711945740146; H = 7.344005520526261; I = 19.02270983468856; K = 9.858118587351978; N = 17.27787165724953; O = 20.28533928839123 }
    15 = @{ B = 10.42018051986992; C = 6.659076658829552; D = 4.821875890675964; E = 20.68409143438207; F = 23.01064399768228; H = 7.344005520526261; I = 19.0296393318344; K = 9.828044578586125; N = 17.28439857916174; O = 20.28781341615436 }
    16 = @{ B = 10.18858976415797; C = 6.497638087373586; D = 4.782390550488811; E = 20.02631573808458; F = 22.97550667843975; H = 7.344005520526261; I = 19.0705659692404; K = 9.654222155501154; N = 17.32228069194875; O = 20.30350577803951 }
    17 = @{ B = 10.04411585700314; C = 6.396311535119256; D = 4.75786468505858; E = 19.61201277992649; F = 22.95571531780267; H = 7.344005520526261; I = 19.09675678467959; K = 9.546318880176347; N = 17.34594850363374; O = 20.314481584174 }
    18 = @{ B = 9.96016027386348; C = 6.33719612439282; D = 4.743648986238585; E = 19.36974492644785; F = 22.94498224192482; H = 7.344005520526261; I = 19.11221855301072; K = 9.483810919519099; N = 17.35971941260347; O = 20.32128949743244 }
    19 = @{ B = 9.931589401009468; C = 6.317037874434779; D = 4.738817277864496; E = 19.28703403892651; F = 22.94146008979354; H = 7.344005520526261; I = 19.11752184425564; K = 9.46257250104499; N = 17.36440914665358; O = 20.32367945823845 }
    20 = @{ B = 10.05958471727972; C = 6.407184492836532; D = 4.760486844689095; E = 19.65652685998972; F = 22.95775487203238; H = 7.344005520526261; I = 19.09392757091869; K = 9.557851904533976; N = 17.34341270071316; O = 20.31326194860139 }
    21 = @{ B = 10.47921190877328; C = 6.7000444966878; D = 4.831974791406905; E = 20.85059236686364; F = 23.02025939722618; H = 7.344005520526261; I = 19.01940849728146; K = 9.87251299306188; N = 17.27475027960656; O = 20.28417968821996 }
    22 = @{ B = 10.74533272070049; C = 6.883885859364369; D = 4.877678147678237; E = 21.5958387922888; F = 23.06685483914202; H = 7.344005520526261; I = 18.9742612710149; K = 10.11412956890204; N = 17.23127247654346; O = 20.26959720436031 }
    23 = @{ B = 10.60406468796979; C = 6.786463254071576; D = 4.853380778848987; E = 21.2012916899622; F = 23.04146281555328; H = 7.344005520526261; I = 18.99803113275984; K = 9.975358314987437; N = 17.25434996333092; O = 20.27697454701332 }
    24 = @{ B = 10.0525940294667; C = 6.402271510404714; D = 4.759301725385036; E = 19.63641475436786; F = 22.95683077914542; H = 7.344005520526261; I = 19.09520539915032; K = 9.55263929221176; N = 17.34455862633701; O = 20.31381179543988 }
    25 = @{ B = 9.42533219184349; C = 5.956043802802212; D = 4.653682617796767; E = 17.7964552231451; F = 22.89042059074011; H = 7.344005520526261; I = 19.2152117968106; K = 9.089243263582203; N = 17.44790188080131; O = 20.37245473735731 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

Write-Output "Updated $($data.Count) rows"